$wb = $excel.ActiveWorkbook

# Sorted genus lists (header excluded) per sheet index (1-based, matches original sheet order)
$sheetData = @{}
$sheetData[1] = @('Acetobacterium', 'Bacillus', 'Clostridium', 'Corynebacterium', 'Desulfobacterium', 'Desulfobulbus', 'Desulfotomaculum', 'Desulfovibrio', 'Gallionella', 'Micrococcus', 'Novosphingobium', 'Propionibacterium', 'Pseudomonas', 'Shewanella', 'Staphylococcus', 'Streptococcus', 'Thiobacillus')
$sheetData[2] = @('Brachybacterium', 'Brevibacterium', 'Bulleidia', 'Enterococcus', 'Gelria', 'Legionella', 'Mycobacterium', 'Mycoplana', 'Neisseria', 'Oerskovia', 'Opitutus', 'Oxobacter', 'Paracoccus', 'Phenylobacterium', 'Porphyrobacter', 'Prevotella', 'Pseudarthrobacter', 'Pseudoalteromonas', 'Tepidimonas')
$sheetData[3] = @('Achromobacter', 'Acidisoma', 'Acidovorax', 'Aestuariimicrobium', 'Afipia', 'Anoxybacillus', 'Beta_proteobacterium', 'Blastomonas', 'Bradyrhizobium', 'Brevundimonas', 'Candidatus_desulforudis', 'Caulobacter', 'Chryseobacterium', 'Clostridium_sensu_stricto_12', 'Cutibacterium', 'Dechloromonas', 'Desulfomicrobium', 'Desulfosporosinus', 'Enhydrobacter', 'Erysipelothrix', 'Flavisolibacter', 'Geothrix', 'Herbaspirillum', 'Hydrogenophaga', 'Methylocystis', 'Nitrospira', 'Oxalobacteraceae_unclassified', 'Phreatobacter', 'Propionivibrio', 'Pseudorhodoferax', 'Pseudoxanthomonas', 'Ralstonia', 'Ruminiclostridium_1', 'Sediminibacterium', 'Silanimonas', 'Simplicispira', 'Smithella', 'Sphingobium', 'Sphingomonas', 'Sphingopyxis', 'Syntrophus', 'Tessaracoccus', 'Thermincola', 'Treponema', 'Variovorax', 'Wchb1-05')
$sheetData[4] = @('Desulfobacterium', 'Desulfobulbus', 'Gallionella', 'Shewanella')
$sheetData[5] = @('Azospira', 'Clostridium', 'Corynebacterium', 'Halomonas', 'Novosphingobium', 'Psb-m-3', 'Streptococcus')
$sheetData[6] = @('Clostridium', 'Corynebacterium', 'Novosphingobium', 'Streptococcus', 'Thiobacillus')
$sheetData[7] = @('Acetobacterium', 'Bacillus', 'Clostridium', 'Corynebacterium', 'Desulfotomaculum', 'Desulfovibrio', 'Micrococcus', 'Novosphingobium', 'Propionibacterium', 'Pseudomonas', 'Staphylococcus', 'Streptococcus')
$sheetData[8] = @('Clostridium', 'Corynebacterium', 'Novosphingobium', 'Streptococcus')

for ($i = 1; $i -le 8; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $names = $sheetData[$i]
    for ($r = 0; $r -lt $names.Length; $r++) {
        $ws.Cells.Item($r + 2, 1).Value = $names[$r]
    }
    $ws.Columns("B").Delete()
}

# New 'components' sheet: union of the composite sheets (pure_checked, checked_core, checked_usual, core_usual)
$components = @('Acetobacterium', 'Azospira', 'Bacillus', 'Brachybacterium', 'Brevibacterium', 'Bulleidia', 'Clostridium', 'Corynebacterium', 'Desulfotomaculum', 'Desulfovibrio', 'Enterococcus', 'Gelria', 'Halomonas', 'Legionella', 'Micrococcus', 'Mycobacterium', 'Mycoplana', 'Neisseria', 'Novosphingobium', 'Oerskovia', 'Opitutus', 'Oxobacter', 'Paracoccus', 'Phenylobacterium', 'Porphyrobacter', 'Prevotella', 'Propionibacterium', 'Psb-m-3', 'Pseudarthrobacter', 'Pseudoalteromonas', 'Pseudomonas', 'Staphylococcus', 'Streptococcus', 'Tepidimonas', 'Thiobacillus')

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newws = $wb.Worksheets.Add($null, $lastSheet)
$newws.Name = "components"

$headerSrc = $wb.Worksheets.Item(1).Range("A1")
$headerSrc.Copy()
$newws.Range("A1").PasteSpecial(-4122)
$newws.Range("A1").Value = "Genus"

for ($r = 0; $r -lt $components.Length; $r++) {
    $newws.Cells.Item($r + 2, 1).Value = $components[$r]
}

# Restore the originally active sheet (adding a worksheet makes it active by default)
$wb.Worksheets.Item(1).Activate()

